# Apply river-results update (May 2024) to the worksheet:
#  - Update recalculated statistic values for existing rows (2014-2018 and
#    2015-2019 .. 2018-2022 groups) whose underlying figures were refreshed.
#  - Append a new "2019 - 2023" results block (rows 82-97).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply updated cell values for existing rows (per diff) ---
$ws.Range("G4").Value = 0.0064757569246722
$ws.Range("G5").Value = 0.0064757569246722
$ws.Range("F7").Value = 0.0013
$ws.Range("G7").Value = 0.0044187328878773
$ws.Range("L7").Value = 0.00165
$ws.Range("M7").Value = 0.00558
$ws.Range("F8").Value = 0.0013
$ws.Range("G8").Value = 0.0044187328878773
$ws.Range("L8").Value = 0.00165
$ws.Range("M8").Value = 0.00558
$ws.Range("F12").Value = 0.5413
$ws.Range("G12").Value = 0.58831724137931
$ws.Range("H12").Value = 1.6235
$ws.Range("N12").Value = 1.09119
$ws.Range("F13").Value = 0.5413
$ws.Range("G13").Value = 0.58831724137931
$ws.Range("H13").Value = 1.6235
$ws.Range("N13").Value = 1.09119
$ws.Range("G16").Value = 0.0195407681857319
$ws.Range("G17").Value = 0.0195407681857319
$ws.Range("G20").Value = 0.0054279956689514
$ws.Range("G21").Value = 0.0054279956689514
$ws.Range("F23").Value = 0.00163
$ws.Range("G23").Value = 0.0047053875790438
$ws.Range("L23").Value = 0.00169
$ws.Range("M23").Value = 0.00753
$ws.Range("N23").Value = 0.0163
$ws.Range("F24").Value = 0.00163
$ws.Range("G24").Value = 0.0047053875790438
$ws.Range("L24").Value = 0.00169
$ws.Range("M24").Value = 0.00753
$ws.Range("N24").Value = 0.0163
$ws.Range("G28").Value = 0.573149152542373
$ws.Range("H28").Value = 1.6235
$ws.Range("M28").Value = 0.91145
$ws.Range("G29").Value = 0.573149152542373
$ws.Range("H29").Value = 1.6235
$ws.Range("M29").Value = 0.91145
$ws.Range("G32").Value = 0.0186841449961432
$ws.Range("G33").Value = 0.0186841449961432
$ws.Range("G36").Value = 0.0047196939159797
$ws.Range("L36").Value = 0.00327
$ws.Range("G37").Value = 0.0047196939159797
$ws.Range("L37").Value = 0.00327
$ws.Range("F39").Value = 0.00169
$ws.Range("G39").Value = 0.0049234045477168
$ws.Range("L39").Value = 0.00263
$ws.Range("M39").Value = 0.00784
$ws.Range("N39").Value = 0.01544
$ws.Range("F40").Value = 0.00169
$ws.Range("G40").Value = 0.0049234045477168
$ws.Range("L40").Value = 0.00263
$ws.Range("M40").Value = 0.00784
$ws.Range("N40").Value = 0.01544
$ws.Range("G44").Value = 0.5061637931034479
$ws.Range("G45").Value = 0.5061637931034479
$ws.Range("G48").Value = 0.0189393235973336
$ws.Range("G49").Value = 0.0189393235973336
$ws.Range("F52").Value = 0.00383
$ws.Range("G52").Value = 0.0044911621503161
$ws.Range("F53").Value = 0.00383
$ws.Range("G53").Value = 0.0044911621503161
$ws.Range("G55").Value = 0.0059385811682667
$ws.Range("L55").Value = 0.00339
$ws.Range("M55").Value = 0.008619999999999999
$ws.Range("N55").Value = 0.0194
$ws.Range("G56").Value = 0.0059385811682667
$ws.Range("L56").Value = 0.00339
$ws.Range("M56").Value = 0.008619999999999999
$ws.Range("N56").Value = 0.0194
$ws.Range("G60").Value = 0.456784482758621
$ws.Range("I60").Value = 0.9559
$ws.Range("G61").Value = 0.456784482758621
$ws.Range("I61").Value = 0.9559
$ws.Range("G64").Value = 0.0229220822180233
$ws.Range("G65").Value = 0.0229220822180233
$ws.Range("G68").Value = 0.0047516279163665
$ws.Range("L68").Value = 0.00301
$ws.Range("G69").Value = 0.0047516279163665
$ws.Range("L69").Value = 0.00301
$ws.Range("F71").Value = 0.00316
$ws.Range("G71").Value = 0.0061479423804353
$ws.Range("L71").Value = 0.00425
$ws.Range("M71").Value = 0.008619999999999999
$ws.Range("N71").Value = 0.01398
$ws.Range("F72").Value = 0.00316
$ws.Range("G72").Value = 0.0061479423804353
$ws.Range("L72").Value = 0.00425
$ws.Range("M72").Value = 0.008619999999999999
$ws.Range("N72").Value = 0.01398
$ws.Range("G76").Value = 0.448439655172414
$ws.Range("H76").Value = 0.9645
$ws.Range("L76").Value = 0.34575
$ws.Range("G77").Value = 0.448439655172414
$ws.Range("H77").Value = 0.9645
$ws.Range("L77").Value = 0.34575

$newRowsData = @{}
$newRowsData[82] = @('Mangatainoka at Hukanui', 'ASPM', 'B', '2019 - 2023', 'RepSite', 0.458, 0.4658, 0.535, 0.535, $null, $null, 0.458, 0.52625, 0.535, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', $null)
$newRowsData[83] = @('Mangatainoka at Hukanui', 'Chlorophyll A (92nd Percentile)', 'A', '2019 - 2023', 'RepSite', 7.75, 9.69875, 39.5, 29.5, $null, $null, 8, 16.6, 25.9, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'mg chl-a /m2')
$newRowsData[84] = @('Mangatainoka at Hukanui', 'DRP (95th Percentile)', 'A', '2019 - 2023', 'RepSite', 0.004, 0.0049264041158838, 0.019, 0.0111, $null, $null, 0.00382, 0.00747, 0.008999999999999999, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'mg/L')
$newRowsData[85] = @('Mangatainoka at Hukanui', 'DRP (Median)', 'A', '2019 - 2023', 'RepSite', 0.004, 0.0049264041158838, 0.019, 0.0111, $null, $null, 0.00382, 0.00747, 0.008999999999999999, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'mg/L')
$newRowsData[86] = @('Mangatainoka at Hukanui', 'MCI', 'C', '2019 - 2023', 'RepSite', 109.57, 107.914, 115, 115, $null, $null, 109.57, 113.6, 115, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', $null)
$newRowsData[87] = @('Mangatainoka at Hukanui', 'Ammoniacal-N (95th Percentile)', 'A', '2019 - 2023', 'RepSite', 0.0031, 0.005996036292652, 0.0501027713572412, 0.02534, $null, $null, 0.00313, 0.00898, 0.0135, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'mg NH4-N/L')
$newRowsData[88] = @('Mangatainoka at Hukanui', 'Ammoniacal-N (Median)', 'A', '2019 - 2023', 'RepSite', 0.0031, 0.005996036292652, 0.0501027713572412, 0.02534, $null, $null, 0.00313, 0.00898, 0.0135, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'mg NH4-N/L')
$newRowsData[89] = @('Mangatainoka at Hukanui', 'Nitrate-N (95th Percentile)', 'A', '2019 - 2023', 'RepSite', 0.379, 0.434745762711864, 0.961, 0.8493000000000001, $null, $null, 0.3365, 0.74999, 0.84478, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'mg NO3-N/L')
$newRowsData[90] = @('Mangatainoka at Hukanui', 'Nitrate-N (Median)', 'A', '2019 - 2023', 'RepSite', 0.379, 0.434745762711864, 0.961, 0.8493000000000001, $null, $null, 0.3365, 0.74999, 0.84478, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'mg NO3-N/L')
$newRowsData[91] = @('Mangatainoka at Hukanui', 'QMCI', 'B', '2019 - 2023', 'RepSite', 6.258, 6.1476, 6.79, 6.79, $null, $null, 6.258, 6.6955, 6.79, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', $null)
$newRowsData[92] = @('Mangatainoka at Hukanui', 'Soluble Inorganic Nitrogen (95th Percentile)', $null, '2019 - 2023', 'RepSite', 0.389, 0.447313559322034, 0.9645, 0.86315, $null, $null, 0.34575, 0.75784, 0.85468, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'g/m3')
$newRowsData[93] = @('Mangatainoka at Hukanui', 'Soluble Inorganic Nitrogen (Median)', $null, '2019 - 2023', 'RepSite', 0.389, 0.447313559322034, 0.9645, 0.86315, $null, $null, 0.34575, 0.75784, 0.85468, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'g/m3')
$newRowsData[94] = @('Mangatainoka at Hukanui', 'Total Nitrogen (95th Percentile)', $null, '2019 - 2023', 'RepSite', 0.55, 0.595932203389831, 1, 0.9255, $null, $null, 0.495, 0.8441, 0.91, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'g/m3')
$newRowsData[95] = @('Mangatainoka at Hukanui', 'Total Nitrogen (Median)', $null, '2019 - 2023', 'RepSite', 0.55, 0.595932203389831, 1, 0.9255, $null, $null, 0.495, 0.8441, 0.91, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'g/m3')
$newRowsData[96] = @('Mangatainoka at Hukanui', 'Total Phosphorus (95th Percentile)', $null, '2019 - 2023', 'RepSite', 0.011, 0.0289491525423729, 0.272, 0.10765, $null, $null, 0.008999999999999999, 0.04247, 0.08892, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'g/m3')
$newRowsData[97] = @('Mangatainoka at Hukanui', 'Total Phosphorus (Median)', $null, '2019 - 2023', 'RepSite', 0.011, 0.0289491525423729, 0.272, 0.10765, $null, $null, 0.008999999999999999, 0.04247, 0.08892, 1829637.96, 5505808.3, 'Tararua District', 'Manawatū', 'Mangatainoka', 'Mana_8b', 'g/m3')

foreach ($rowNum in $newRowsData.Keys) {
    $rowVals = $newRowsData[$rowNum]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $cellVal = $rowVals[$i]
        if ($null -ne $cellVal) {
            $ws.Cells.Item($rowNum, $i + 1).Value = $cellVal
        }
    }
}
